# Update scripts with new TPM values.
# Root change: recomputed "Receptor average/total expression value" for
# Eng receptor in ECs (rows 2 and 5) using new TPM values, which in turn
# changes the derived specificity / edge-weight columns (O:T) for every
# row of the sheet (the M,N,O,P columns are normalized per sending-cluster
# group of 3 rows, and Q,R,S,T are edge weights normalized across all rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> Eng in ECs)
$ws.Range("M2").Value = 201.098592
$ws.Range("N2").Value = 603.295776
$ws.Range("O2").Value = 0.7918622805845071
$ws.Range("P2").Value = 0.791862280584507
$ws.Range("Q2").Value = 860.3282655432
$ws.Range("R2").Value = 7742.954389888801
$ws.Range("S2").Value = 0.7651943461005664
$ws.Range("T2").Value = 0.7651943461005661

# Row 3 (FAPs -> Eng in FAPs)
$ws.Range("O3").Value = 0.1414593902976603
$ws.Range("P3").Value = 0.1414593902976603
$ws.Range("S3").Value = 0.1366953930154414
$ws.Range("T3").Value = 0.1366953930154414

# Row 4 (FAPs -> Eng in MuSCs)
$ws.Range("O4").Value = 0.0666783291178327
$ws.Range("P4").Value = 0.06667832911783268
$ws.Range("S4").Value = 0.06443277031801152
$ws.Range("T4").Value = 0.06443277031801151

# Row 5 (MuSCs -> Eng in ECs)
$ws.Range("M5").Value = 201.098592
$ws.Range("N5").Value = 603.295776
$ws.Range("O5").Value = 0.7918622805845071
$ws.Range("P5").Value = 0.791862280584507
$ws.Range("Q5").Value = 29.98346490288001
$ws.Range("R5").Value = 269.85118412592
$ws.Range("S5").Value = 0.02666793448394087
$ws.Range("T5").Value = 0.02666793448394087

# Row 6 (MuSCs -> Eng in FAPs)
$ws.Range("O6").Value = 0.1414593902976603
$ws.Range("P6").Value = 0.1414593902976603
$ws.Range("S6").Value = 0.004763997282218867
$ws.Range("T6").Value = 0.004763997282218865

# Row 7 (MuSCs -> Eng in MuSCs)
$ws.Range("O7").Value = 0.0666783291178327
$ws.Range("P7").Value = 0.06667832911783268
$ws.Range("S7").Value = 0.002245558799821181
$ws.Range("T7").Value = 0.002245558799821181
